# Auto-generated script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.180.46'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '1.854.33'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +1.03%  '
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").Value = '''310.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").Value = '''0.4777'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.76%  '
$ws.Range("D8").Value = '''0.3700'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.02%  '
$ws.Range("D9").Value = '''0.07278'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.63%  '
$ws.Range("D10").Value = '''0.9343'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.72%  '
$ws.Range("D11").Value = '''19.95'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.81%  '
$ws.Range("D12").Value = '''0.07818'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").Value = '1.855.95'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("E14").Value = '  +2.06%  '
$ws.Range("D15").Value = '''6.512'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.64%  '
$ws.Range("D16").Value = '''89.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.53%  '
$ws.Range("D17").Value = '''1.019'
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("D20").Value = '27.204.22'
$ws.Range("E20").Value = '  +0.90%  '
$ws.Range("D21").Value = '''14.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("D22").Value = '''5.084'
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("D24").Value = '''1.946'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").Value = '''153.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("E26").Value = '  +1.01%  '
$ws.Range("D27").Value = '''1.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.17%  '
$ws.Range("D28").Value = '''115.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.71%  '
$ws.Range("D29").Value = '''4.935'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.08%  '
$ws.Range("D30").Value = '''0.08877'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").Value = '''3.310'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.70%  '
$ws.Range("D32").Value = '''1.182'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("D33").Value = '''4.551'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.62%  '
$ws.Range("D34").Value = '''0.7380'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("D35").Value = '''2.691'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.43%  '
$ws.Range("D36").Value = '''1.117'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.91%  '
$ws.Range("D37").Value = '''0.01998'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.95%  '
$ws.Range("D38").Value = '''0.05257'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.16%  '
$ws.Range("D39").Value = '''2.978'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.66%  '
$ws.Range("D40").Value = '''0.5303'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.99%  '
$ws.Range("D41").Value = '''7.052'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = '''0.1529'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("D43").Value = '''8.330'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.33%  '
$ws.Range("D44").Value = '''10.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.67%  '
$ws.Range("D45").Value = '''0.4754'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.21%  '
$ws.Range("D46").Value = '''1.018'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.98%  '
$ws.Range("D47").Value = '''102.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.84%  '
$ws.Range("D48").Value = '''1.628'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.18%  '
$ws.Range("D49").Value = '''66.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.83%  '
$ws.Range("D50").Value = '''0.06066'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = '''0.8945'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.21%  '
